$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.910.49"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.815.50"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.29"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4644"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3662"
$ws.Range("E8").Value = "  -0.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07354"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  -0.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.30"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.841.87"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.375"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07087"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.46"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008713"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.64"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").Value = "26.941.34"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.300"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "2.049.25"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.895"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.85"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.31"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.131"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.256"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.28"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08900"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7574"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.157"
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.482"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.911"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05276"
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01945"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.982"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.230"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5293"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.285"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1655"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.423"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4868"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.45"
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.28"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.660"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06292"
